# Revert "Drop in files from RMI script"
# Re-add the "Texas Notes" worksheet (with its notes content) between the
# "Calculations" and "PEUDfSbQL" sheets, and restore the selection/active
# sheet state that came with it.

$wb = $excel.ActiveWorkbook

# --- Insert the "Texas Notes" sheet just before "PEUDfSbQL" -----------------
$peSheet = $wb.Worksheets.Item("PEUDfSbQL")
$notes = $wb.Worksheets.Add($peSheet)
$notes.Name = "Texas Notes"

# --- Populate the notes content ---------------------------------------------
$notes.Range("A1").Value = "They are just comparing the efficiency of new appliances:"
$notes.Range("A2").Value = "standard versus energy star rebate qualifying"
$notes.Range("A3").Value = "i.e., the point of this spreadsheet is to estimate how much a household's energy consumption would change"
$notes.Range("A4").Value = "if they decide to use a rebate to get a higher-efficiency appliance instead of just opting for the cheap alternative. "
$notes.Range("A6").Value = "I think the method makes sense"
$notes.Range("A8").Value = "And there's no reason that Texas should be different. "
$notes.Range("A9").Value = "New technology in Texas should be as efficient as new technology across the US. "
$notes.Range("A10").Value = "The only difference might be if Texas rebates incentivize a different level of efficiency than"
$notes.Range("A11").Value = "national rebates do, but some of the other sources used in the building input files seem"
$notes.Range("A12").Value = "to indicate that Texas doesn't usually have appliance rebates on top of the national ones. "
$notes.Range("A14").Value = "So it's a good assumption that if a Texan uses a rebate to buy a more efficient appliance, that"
$notes.Range("A15").Value = "rebate will be a national one and it will be based on national energy star standards."

# --- Restore view/selection state on the other sheets -----------------------
$dataSheet = $wb.Worksheets.Item("Data")
[void]$dataSheet.Activate()
[void]$dataSheet.Range("C6").Select()

$calcSheet = $wb.Worksheets.Item("Calculations")
[void]$calcSheet.Activate()
[void]$calcSheet.Range("B18").Select()

# --- Make "PEUDfSbQL" the active/selected sheet again ------------------------
# (re-fetch by name: the sheet collection shifted when "Texas Notes" was
# inserted, so the old $peSheet handle no longer points at PEUDfSbQL)
$peSheet = $wb.Worksheets.Item("PEUDfSbQL")
[void]$peSheet.Activate()
[void]$peSheet.Range("B7").Select()
